# Adds season record columns (Wins, Losses, Ties) to the NYY 2018 worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): AD1 = Wins, AE1 = Losses, AF1 = Ties ---
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header formatting (style) from an existing header cell (A1) onto
# the three new header cells so they match the rest of the header row.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122, $null)

# --- Data rows (2-51): every team row gets the same season record ---
$wins = 100
$losses = 62
$ties = 0

for ($r = 2; $r -le 51; $r++) {
    $ws.Cells.Item($r, 30).Value = $wins    # column AD
    $ws.Cells.Item($r, 31).Value = $losses  # column AE
    $ws.Cells.Item($r, 32).Value = $ties    # column AF
}
